$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.879.70'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('D3').Value = '1.904.51'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5005'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2997'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06860'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.82%  '
$ws.Range('D10').Value = '1.906.30'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.44'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07358'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '91.71'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.125'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6825'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.91%  '
$ws.Range('D16').Value = '30.862.59'
$ws.Range('E16').Value = '  +2.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008068'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').Value = '2.150.76'
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9991'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.882'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '184.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +36.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.104'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.392'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.20'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E27').Value = '  +11.08%  '
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.397'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09003'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.087'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05290'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7480'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.70%  '
$ws.Range('E35').Value = '  +3.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.672'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01929'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +17.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.726'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.192'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9443'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4401'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '106.38'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.868'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.800'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1360'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05857'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.3935'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.598'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.48'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.35%  '
$ws.Range('E51').Value = '  +4.20%  '
